# Update countries & provincias Spain
#
# This script applies the data refresh that happened between the
# "08:52" and "09:22" snapshots of the COVID-19 "paises" workbook:
#   - Re-labels a handful of rows whose country changed because new
#     countries (Finlandia, Armenia) were inserted ahead of their
#     former neighbours in the source feed, shifting the rows below
#     them down by one.
#   - Refreshes the numeric statistics (Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#     for every row whose figures moved.
#   - Updates the "Datos actualizados a ..." timestamp banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Timestamp banner (row 1)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 09:22"

# ---------------------------------------------------------------
# 2. Country label shifts caused by newly inserted rows upstream
# ---------------------------------------------------------------
$ws.Range("A43").Value = "Finlandia"
$ws.Range("A44").Value = "Serbia"
$ws.Range("A45").Value = "Tailandia"
$ws.Range("A46").Value = "Emiratos Arabes Unidos"

$ws.Range("A69").Value = "Armenia"
$ws.Range("A70").Value = "Bielorrusia"

# ---------------------------------------------------------------
# 3. Updated statistics per row
#    Columns: B=Casos totales, C=Nuevos casos, D=Casos activos,
#             E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes
# ---------------------------------------------------------------

function Set-Row($rowNum, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($rowNum, 2).Value = $b
    $ws.Cells.Item($rowNum, 3).Value = $c
    $ws.Cells.Item($rowNum, 4).Value = $d
    $ws.Cells.Item($rowNum, 5).Value = $e
    $ws.Cells.Item($rowNum, 6).Value = $f
    $ws.Cells.Item($rowNum, 7).Value = $g
    $ws.Cells.Item($rowNum, 8).Value = $h
}

# Row 4: Estados Unidos
Set-Row 4 400546 211 21711 365978 9169 16 12857

# Row 13: Suiza
Set-Row 13 22328 75 8704 12800 391 3 824

# Row 18: Austria
Set-Row 18 12709 70 4046 8420 243 0 243

# Row 27: India
Set-Row 27 5360 9 468 4728 0 4 164

# Row 32: Rumania
Set-Row 32 4417 0 460 3752 274 8 205

# Row 43: now Finlandia
Set-Row 43 2487 179 300 2153 81 0 34

# Row 44: now Serbia
Set-Row 44 2447 0 118 2268 109 0 61

# Row 45: now Tailandia
Set-Row 45 2369 111 888 1451 61 3 30

# Row 46: now Emiratos Arabes Unidos
Set-Row 46 2359 0 186 2161 1 0 12

# Row 69: now Armenia
Set-Row 69 881 28 114 758 30 1 9

# Row 70: now Bielorrusia
Set-Row 70 861 0 54 794 31 0 13

# Row 85: Uzbekistan
Set-Row 85 534 14 30 501 8 1 3

# Row 163: Antigua y Barbuda
Set-Row 163 19 0 0 17 1 1 2
